$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 235
$ws.Range("I19").Value = 280.4
$ws.Range("J19").Value = 197.16667
$ws.Range("K19").Value = 280.4
$ws.Range("L19").Value = 197.16667
$ws.Range("M19").Value = -105.4
$ws.Range("N19").Value = -547.1666700000001

$ws.Range("H40").Value = 4999.1665
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825

$ws.Range("H95").Value = 33999
$ws.Range("J95").Value = 33999
$ws.Range("L95").Value = 33999
$ws.Range("N95").Value = -39491

$ws.Range("H100").Value = 3986.1428
$ws.Range("I100").Value = 3227
$ws.Range("K100").Value = 3227
$ws.Range("M100").Value = -2686

$ws.Range("H113").Value = 2844.5454
$ws.Range("I113").Value = 2361.375
$ws.Range("J113").Value = 4133
$ws.Range("K113").Value = 2361.375
$ws.Range("L113").Value = 4133
$ws.Range("M113").Value = 892.625
$ws.Range("N113").Value = -10641

$ws.Range("H138").Value = 3668.3872
$ws.Range("I138").Value = 1483.1666
$ws.Range("J138").Value = 4192.84
$ws.Range("K138").Value = 4449.4998
$ws.Range("L138").Value = 12578.52
$ws.Range("M138").Value = 690.5002000000004
$ws.Range("N138").Value = -22858.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 47621844
$ws.Range("I97").Value = 47621844
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 47621844
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -47621348
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 4280.5
$ws.Range("J122").Value = 2520.2
$ws.Range("L122").Value = 7560.599999999999
$ws.Range("N122").Value = -12460.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1553.3684
$ws.Range("I86").Value = 1153.7646
$ws.Range("J86").Value = 4950
$ws.Range("K86").Value = 1153.7646
$ws.Range("L86").Value = 4950
$ws.Range("M86").Value = -30.76459999999997
$ws.Range("N86").Value = -7196

$ws.Range("H89").Value = 1553.3684
$ws.Range("I89").Value = 1153.7646
$ws.Range("J89").Value = 4950
$ws.Range("K89").Value = 5768.823
$ws.Range("L89").Value = 24750
$ws.Range("M89").Value = -152.8230000000003
$ws.Range("N89").Value = -35982

$ws.Range("H94").Value = 158095.14
$ws.Range("I94").Value = 220971.2
$ws.Range("K94").Value = 220971.2
$ws.Range("M94").Value = -220520.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2604.25
$ws.Range("I105").Value = 2547.7144
$ws.Range("K105").Value = 2547.7144
$ws.Range("M105").Value = -800.7143999999998

$ws.Range("H122").Value = 2067
$ws.Range("I122").Value = 1910.3334
$ws.Range("K122").Value = 5731.0002
$ws.Range("M122").Value = -3281.0002

$ws.Range("H141").Value = 629628.7
$ws.Range("J141").Value = 629628.7
$ws.Range("L141").Value = 629628.7
$ws.Range("N141").Value = -639988.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 6693.75
$ws.Range("I99").Value = 6693.75
$ws.Range("K99").Value = 20081.25
$ws.Range("M99").Value = -17835.25

$ws.Range("H131").Value = 1010.25
$ws.Range("I131").Value = 726
$ws.Range("K131").Value = 2178
$ws.Range("M131").Value = 2862

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H9").Value = 353.5
$ws.Range("I9").Value = 353.5
$ws.Range("K9").Value = 353.5
$ws.Range("M9").Value = -183.5

$ws.Range("H70").Value = 250000600
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 500000000
$ws.Range("K70").Value = 1200
$ws.Range("L70").Value = 500000000
$ws.Range("M70").Value = -930
$ws.Range("N70").Value = -500000540

$ws.Range("H73").Value = 250000600
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 500000000
$ws.Range("K73").Value = 1200
$ws.Range("L73").Value = 500000000
$ws.Range("M73").Value = -264
$ws.Range("N73").Value = -500001872

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 513.3333
$ws.Range("I9").Value = 270
$ws.Range("K9").Value = 270
$ws.Range("M9").Value = -46

$ws.Range("H22").Value = 294.69565
$ws.Range("I22").Value = 182.85715
$ws.Range("K22").Value = 182.85715
$ws.Range("M22").Value = 112.14285

$ws.Range("H27").Value = 294.69565
$ws.Range("I27").Value = 182.85715
$ws.Range("K27").Value = 182.85715
$ws.Range("M27").Value = -75.85714999999999

$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9705

$ws.Range("H31").Value = 4151.75
$ws.Range("I31").Value = 2071.6667
$ws.Range("J31").Value = 5399.8
$ws.Range("K31").Value = 2071.6667
$ws.Range("L31").Value = 5399.8
$ws.Range("M31").Value = -1823.6667
$ws.Range("N31").Value = -5895.8

$ws.Range("H46").Value = 513.5714
$ws.Range("I46").Value = 239
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 239
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = -51
$ws.Range("N46").Value = -1576

$ws.Range("H68").Value = 6899.091
$ws.Range("J68").Value = 7222.222
$ws.Range("L68").Value = 7222.222
$ws.Range("N68").Value = -8720.222

$ws.Range("H71").Value = 6899.091
$ws.Range("J71").Value = 7222.222
$ws.Range("L71").Value = 36111.11
$ws.Range("N71").Value = -43599.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5999.5
$ws.Range("I62").Value = 5999.5
$ws.Range("K62").Value = 5999.5
$ws.Range("M62").Value = -5375.5

$ws.Range("H65").Value = 5999.5
$ws.Range("I65").Value = 5999.5
$ws.Range("K65").Value = 29997.5
$ws.Range("M65").Value = -26877.5

$ws.Range("H96").Value = 3414.7144
$ws.Range("I96").Value = 3200.6
$ws.Range("J96").Value = 3950
$ws.Range("K96").Value = 3200.6
$ws.Range("L96").Value = 3950
$ws.Range("M96").Value = -1827.6
